# Update "想去人数" (interest count) values on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 1496
    $ws.Range("F6").Value = 36
    $ws.Range("F9").Value = 292
}
